$d = $word.ActiveDocument

# --- 1. Course title: "CSIS 3540 - Client Server Systems" -> "CSIS 3280 - Web Scripting"
$d.Content.Find.Execute("Course:   CSIS 3540", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "Course:   CSIS 3280", 2) | Out-Null
$d.Content.Find.Execute("Client Server Systems", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "Web Scripting", 2) | Out-Null

# --- 2. Due Date day: "Due Date: 3/04/2019" -> "Due Date: 5/04/2019"
$d.Content.Find.Execute("Due Date: 3", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "Due Date: 5", 2) | Out-Null

# --- 3. (Rebouças) proofErr spell markers only - no visible text change, not
#        reachable through the Word object model (no scriptable "run spell
#        check" verb); skipped.

# --- 4. Project-description paragraph rewrite.
$d.Content.Find.Execute("project but we also", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "project, but we also", 2) | Out-Null

$d.Content.Find.Execute("our main entities will also contain CRUD operations", `
                         $false, $true, $false, $false, $false, `
                         $true, 1, $false, "all entities will also contain CRUD operations", 2) | Out-Null

$d.Content.Find.Execute("on to be added to the database.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "on to be added to the database. If time permits, we will add a REST API to one of our tables for additional simplicity of data retrieval from the database, while adding another layer of complexity to the project.", `
                         2) | Out-Null

# --- 5. "Facilities_ship:" -> "Facilities_Ship:" (capitalize the S) and move
#        the stray "_GoBack" bookmark here (it was left near "...and log
#        file." from the previous save; Word re-drops it at the most recent
#        edit point, which is this heading after this edit).
$d.Content.Find.Execute("Facilities_ship:", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "Facilities_Ship:", 2) | Out-Null

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$afterS = $d.Content
$afterS.Find.Execute("Facilities_S") | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($afterS.End, $afterS.End)) | Out-Null

# --- 6. "One way cruise" -> "One-way cruise"
$d.Content.Find.Execute("One way cruise", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "One-way cruise", 2) | Out-Null

# --- 7 & 9. "Facilities_Ship and components" / "Lindsey: ... Facilities_Ship.
#        ... Facilities_Ships Page." only gain proofErr spell-check markers
#        around "Facilities_Ship"/"Facilities_Ships" - no visible text
#        change, so nothing to do here (proofErr is not scriptable via the
#        object model).

# --- 10. "To add a new facility, choose Ship and Facility to add to that
#        ship from dropdown boxes" -> "To add a new facility to a ship,
#        choose Ship and Facility from dropdown boxes"
$d.Content.Find.Execute( `
    "To add a new facility, choose Ship and Facility to add to that ship from dropdown boxes", `
    $false, $false, $false, $false, $false, `
    $true, 1, $false, `
    "To add a new facility to a ship, choose Ship and Facility from dropdown boxes", `
    2) | Out-Null

# --- 11. "edit will be updated to the list" -> "information will be updated
#        to the list"
$d.Content.Find.Execute("edit will be updated to the list", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "information will be updated to the list", 2) | Out-Null
